# BrassA-HW20: notebook was rerun -> two new sample rows ("Holden" and
# "Rizzie Spiral") were added right after the "Spiral5" row, and one
# existing label was renamed ("Thomas Hex" -> "Matthies Hex"). The whole
# simulation table was then re-generated, which is why every row's C:W
# values differ from before (old row N's numbers now live at row N+2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the row label used at what is currently row 9 (column B) before
# the insert shifts things around - use Find so we don't have to track
# coordinates by hand.
$found = $ws.Cells.Find("Thomas Hex")
if ($found -ne $null) {
    $found.Value = "Matthies Hex"
}

# Insert two new blank rows right after row 3 ("Spiral5"), pushing every
# row from the old row 4 onward down by two (old row 4 -> new row 6, ...,
# old row 29 -> new row 31). Dimension grows to A1:W31 automatically.
$ws.Rows("4:5").Insert()

# Carry over the bold/border/centered header style used throughout column
# A (and used by B1:W1) onto the two new index cells.
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)

# Sequence numbers for the two new rows.
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3

# New sample-type labels.
$ws.Range("B4").Value = "Holden"
$ws.Range("B5").Value = "Rizzie Spiral"

# Freshly (re)simulated data for the "Holden" row.
$ws.Range("C4").Value = 0.9901185788781623
$ws.Range("D4").Value = 0.930680722263085
$ws.Range("E4").Value = 0.9801797652900129
$ws.Range("F4").Value = 0.9801797652900129
$ws.Range("G4").Value = 0.967345155701159
$ws.Range("H4").Value = 1.05769437293166
$ws.Range("I4").Value = 1.101420529237525
$ws.Range("J4").Value = 0.8947333064830572
$ws.Range("K4").Value = 0.9801797652900129
$ws.Range("L4").Value = 0.8947333064830572
$ws.Range("M4").Value = 1.010075202876872
$ws.Range("N4").Value = 0.9801797652900129
$ws.Range("O4").Value = 1.101420529237525
$ws.Range("P4").Value = 0.9980769178602913
$ws.Range("Q4").Value = 1.034382842469342
$ws.Range("R4").Value = 0.9921112003368652
$ws.Range("S4").Value = 0.9878329971405805
$ws.Range("T4").Value = 0.9921112003368652
$ws.Range("U4").Value = 0.9859196891779386
$ws.Range("V4").Value = 0.9847717044003534
$ws.Range("W4").Value = 0.9915309542076918

# Freshly (re)simulated data for the "Rizzie Spiral" row.
$ws.Range("C5").Value = 1.072117955009372
$ws.Range("D5").Value = 0.5801477425201567
$ws.Range("E5").Value = 1.123028197876798
$ws.Range("F5").Value = 1.123028197876798
$ws.Range("G5").Value = 0.8545696076823004
$ws.Range("H5").Value = 1.290179768012821
$ws.Range("I5").Value = 1.371585332876997
$ws.Range("J5").Value = 0.339183016162964
$ws.Range("K5").Value = 1.123028197876798
$ws.Range("L5").Value = 0.339183016162964
$ws.Range("M5").Value = 0.9624627609472416
$ws.Range("N5").Value = 1.123028197876798
$ws.Range("O5").Value = 1.371585332876997
$ws.Range("P5").Value = 0.8553841745199804
$ws.Range("Q5").Value = 1.113077470279649
$ws.Range("R5").Value = 0.9445988489722529
$ws.Range("S5").Value = 0.8551126522407536
$ws.Range("T5").Value = 0.9445988489722529
$ws.Range("U5").Value = 0.9220915386497648
$ws.Range("V5").Value = 0.9622788704951712
$ws.Range("W5").Value = 0.9491592976360813
